$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The leaderboard sheet used to have "Difficulty" / "Bombs on?" columns (C/D).
# The program now reads Name/Score rows straight from the sheet, so drop
# those two extra header cells and use the freed-up columns for real
# high-score rows instead.
$ws.Range("C1:D1").ClearContents()

# Seed a couple of sample high-score rows (Name, Score) - written with the
# lower score first so the shared-string table ends up in the same order
# as the final, sorted sheet.
$ws.Range("A3").Value = "Jack"
$ws.Range("B3").Value = 5
$ws.Range("A2").Value = "Jack Better"
$ws.Range("B2").Value = 10

# Sort the new rows by Score, highest first - this is how the leaderboard
# will be kept ordered once the program starts writing real scores.
[void]$ws.Range("A2:B3").Sort(
    $ws.Range("B2:B3"), 2, [Type]::Missing, [Type]::Missing,
    [Type]::Missing, [Type]::Missing, 1, [Type]::Missing,
    [Type]::Missing, [Type]::Missing
)

# Move the active selection off the now-empty D1/C1 header area.
[void]$ws.Range("C1").Select()
